$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells whose shared-string content changed following the re-processing
# of the data with the newly curated dimensions.
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Comunidad"

# The "aragon" mapping file reference is no longer needed.
$ws.Range("E5").Clear()
